$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "38.796.33"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "2.104.36"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "229.49"
$ws.Range("E5").Value = "  +0.42%  "
Set-TextValue $ws.Range("D6") "0.619"
$ws.Range("E6").Value = "  +0.92%  "
Set-TextValue $ws.Range("D7") "62.50"
$ws.Range("E7").Value = "  +2.33%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.05%  "
Set-TextValue $ws.Range("D10") "0.0845"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  -0.03%  "
Set-TextValue $ws.Range("D12") "15.85"
$ws.Range("E12").Value = "  +7.47%  "
$ws.Range("D13").Value = "2.416.55"
$ws.Range("E13").Value = "  +0.55%  "
Set-TextValue $ws.Range("D14") "22.08"
$ws.Range("E14").Value = "  -1.17%  "
Set-TextValue $ws.Range("D15") "0.808"
$ws.Range("E15").Value = "  +3.82%  "
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "2.097.39"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "38.828.09"
$ws.Range("E18").Value = "  +1.77%  "
Set-TextValue $ws.Range("D19") "72.07"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "0.0₃0842"
$ws.Range("E21").Value = "  +0.55%  "
Set-TextValue $ws.Range("D22") "228.04"
$ws.Range("E22").Value = "  +1.92%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("E25").Value = "  +1.16%  "
Set-TextValue $ws.Range("D26") "172.02"
$ws.Range("E26").Value = "  +1.22%  "
Set-TextValue $ws.Range("D27") "9.59"
$ws.Range("E27").Value = "  +1.49%  "
Set-TextValue $ws.Range("D28") "0.138"
$ws.Range("E28").Value = "  +5.90%  "
$ws.Range("E29").Value = "  +3.41%  "
Set-TextValue $ws.Range("D30") "19.37"
$ws.Range("E30").Value = "  +1.84%  "
Set-TextValue $ws.Range("D31") "2.48"
$ws.Range("E31").Value = "  +3.90%  "
$ws.Range("E32").Value = "  +1.07%  "
Set-TextValue $ws.Range("D33") "4.54"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("E35").Value = "  +2.62%  "
$ws.Range("E36").Value = "  +2.67%  "
$ws.Range("E37").Value = "  +1.58%  "
Set-TextValue $ws.Range("D38") "3.57"
$ws.Range("E38").Value = "  +1.45%  "
$ws.Range("E39").Value = "  -0.01%  "
Set-TextValue $ws.Range("D40") "18.35"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("E41").Value = "  +4.35%  "
Set-TextValue $ws.Range("D42") "101.78"
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("D43").Value = "1.534.89"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("E44").Value = "  -0.91%  "
Set-TextValue $ws.Range("D45") "7.78"
$ws.Range("E45").Value = "  +4.35%  "
Set-TextValue $ws.Range("D46") "0.0912"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("E47").Value = "  +2.04%  "
Set-TextValue $ws.Range("D48") "4.12"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").Value = "2.299.02"
$ws.Range("E51").Value = "  +0.35%  "
